$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 57.2
$ws.Range("I11").Value = 57.2
$ws.Range("K11").Value = 57.2
$ws.Range("M11").Value = 82.8

$ws.Range("H19").Value = 879.4375
$ws.Range("J19").Value = 827.3
$ws.Range("L19").Value = 827.3
$ws.Range("N19").Value = -1177.3

$ws.Range("H96").Value = 125491.5
$ws.Range("J96").Value = 885
$ws.Range("L96").Value = 2655
$ws.Range("N96").Value = -5401

$ws.Range("H100").Value = 2558.4
$ws.Range("I100").Value = 2291.2
$ws.Range("J100").Value = 2825.6
$ws.Range("K100").Value = 2291.2
$ws.Range("L100").Value = 2825.6
$ws.Range("M100").Value = -1750.2
$ws.Range("N100").Value = -3907.6

$ws.Range("H106").Value = 372806.72
$ws.Range("I106").Value = 428730.94
$ws.Range("K106").Value = 428730.94
$ws.Range("M106").Value = -428099.94

$ws.Range("H107").Value = 1178
$ws.Range("I107").Value = 246
$ws.Range("J107").Value = 3042
$ws.Range("K107").Value = 246
$ws.Range("L107").Value = 3042
$ws.Range("M107").Value = 1674
$ws.Range("N107").Value = -6882

$ws.Range("H109").Value = 88141.664
$ws.Range("J109").Value = 88141.664
$ws.Range("L109").Value = 88141.664
$ws.Range("N109").Value = -90915.664

$ws.Range("H111").Value = 1008.3333
$ws.Range("I111").Value = 910
$ws.Range("K111").Value = 2730
$ws.Range("M111").Value = 337

$ws.Range("H137").Value = 559368.4399999999
$ws.Range("I137").Value = 1375.75
$ws.Range("K137").Value = 4127.25
$ws.Range("M137").Value = -1577.25

$ws.Range("H138").Value = 2190.717
$ws.Range("I138").Value = 891.9524
$ws.Range("J138").Value = 3043.0312
$ws.Range("K138").Value = 2675.8572
$ws.Range("L138").Value = 9129.0936
$ws.Range("M138").Value = 2464.1428
$ws.Range("N138").Value = -19409.0936

$ws.Range("H141").Value = 3044.5557
$ws.Range("I141").Value = 3041.2942
$ws.Range("J141").Value = 3100
$ws.Range("K141").Value = 9123.882599999999
$ws.Range("L141").Value = 9300
$ws.Range("M141").Value = -3943.882599999999
$ws.Range("N141").Value = -19660

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8409.705
$ws.Range("I32").Value = 4172.961
$ws.Range("J32").Value = 30017.1
$ws.Range("K32").Value = 4172.961
$ws.Range("L32").Value = 30017.1
$ws.Range("M32").Value = -3885.961
$ws.Range("N32").Value = -30591.1

$ws.Range("H45").Value = 2901.375
$ws.Range("I45").Value = 2535.1667
$ws.Range("K45").Value = 2535.1667
$ws.Range("M45").Value = -2158.1667

$ws.Range("H61").Value = 2539.4614
$ws.Range("I61").Value = 2157
$ws.Range("J61").Value = 3400
$ws.Range("K61").Value = 2157
$ws.Range("L61").Value = 3400
$ws.Range("M61").Value = -1945
$ws.Range("N61").Value = -3824

$ws.Range("H80").Value = 50000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 50000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 50000
$ws.Range("N80").Value = -51996
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 50000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 50000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 150000
$ws.Range("N83").Value = -159984
$ws.Range("M83").ClearContents()

$ws.Range("H132").Value = 2019.4117
$ws.Range("I132").Value = 1518
$ws.Range("J132").Value = 2938.6667
$ws.Range("K132").Value = 4554
$ws.Range("L132").Value = 8816.000100000001
$ws.Range("M132").Value = -2024
$ws.Range("N132").Value = -13876.0001

$ws.Range("H136").Value = 2539.4614
$ws.Range("I136").Value = 2157
$ws.Range("J136").Value = 3400
$ws.Range("K136").Value = 6471
$ws.Range("L136").Value = 10200
$ws.Range("M136").Value = -3921
$ws.Range("N136").Value = -15300

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 41004.707
$ws.Range("J132").Value = 41004.707
$ws.Range("L132").Value = 41004.707
$ws.Range("N132").Value = -51124.707

$ws.Range("H134").Value = 1540.6818
$ws.Range("I134").Value = 1173.421
$ws.Range("K134").Value = 3520.263
$ws.Range("M134").Value = -985.2629999999999

$ws.Range("H135").Value = 44299.855
$ws.Range("J135").Value = 44299.855
$ws.Range("L135").Value = 44299.855
$ws.Range("N135").Value = -54439.855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10639.143
$ws.Range("I31").Value = 1916.1538
$ws.Range("J31").Value = 24814
$ws.Range("K31").Value = 1916.1538
$ws.Range("L31").Value = 24814
$ws.Range("M31").Value = -1621.1538
$ws.Range("N31").Value = -25404

$ws.Range("H34").Value = 10639.143
$ws.Range("I34").Value = 1916.1538
$ws.Range("J34").Value = 24814
$ws.Range("K34").Value = 1916.1538
$ws.Range("L34").Value = 24814
$ws.Range("M34").Value = -1714.1538
$ws.Range("N34").Value = -25218

$ws.Range("H134").Value = 1704.5
$ws.Range("I134").Value = 1748.1428
$ws.Range("J134").Value = 1399
$ws.Range("K134").Value = 5244.428400000001
$ws.Range("L134").Value = 4197
$ws.Range("M134").Value = -2709.428400000001
$ws.Range("N134").Value = -9267

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 4691.6665
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 4691.6665
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 14074.9995
$ws.Range("N54").Value = -15192.9995
$ws.Range("M54").ClearContents()

$ws.Range("H121").Value = 1848.5769
$ws.Range("J121").Value = 1917.8096
$ws.Range("L121").Value = 5753.4288
$ws.Range("N121").Value = -8373.4288

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 18916.062
$ws.Range("J122").Value = 14667.846
$ws.Range("L122").Value = 44003.538
$ws.Range("N122").Value = -48903.538

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 23529.182
$ws.Range("I7").Value = 12928.263
$ws.Range("K7").Value = 12928.263
$ws.Range("M7").Value = -12816.263

$ws.Range("H126").Value = 23529.182
$ws.Range("I126").Value = 12928.263
$ws.Range("K126").Value = 38784.789
$ws.Range("M126").Value = -36314.789

$ws.Range("H132").Value = 8958.375
$ws.Range("I132").Value = 12531
$ws.Range("K132").Value = 37593
$ws.Range("M132").Value = -35063

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 27889.334
$ws.Range("I61").Value = 24429.143
$ws.Range("J61").Value = 40000
$ws.Range("K61").Value = 24429.143
$ws.Range("L61").Value = 40000
$ws.Range("M61").Value = -24137.143
$ws.Range("N61").Value = -40584

$ws.Range("H81").Value = 1917.4762
$ws.Range("J81").Value = 2559.8
$ws.Range("L81").Value = 5119.6
$ws.Range("N81").Value = -7241.6

$ws.Range("H84").Value = 1917.4762
$ws.Range("J84").Value = 2559.8
$ws.Range("L84").Value = 25598
$ws.Range("N84").Value = -36206

$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws.Range("H122").Value = 3494.7273
$ws.Range("I122").Value = 2649.4
$ws.Range("K122").Value = 7948.200000000001
$ws.Range("M122").Value = -5498.200000000001

$ws.Range("H132").Value = 1977779.4
$ws.Range("I132").Value = 1232
$ws.Range("J132").Value = 3954326.8
$ws.Range("K132").Value = 3696
$ws.Range("L132").Value = 11862980.4
$ws.Range("M132").Value = -1166
$ws.Range("N132").Value = -11868040.4
